$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 59
$srcRow = 58

# Seed the new row by copying the previous data row's cell formatting
# (style index) column by column, then overwrite each cell's value with
# the new log entry's data.
for ($col = 1; $col -le 8; $col++) {
    $ws.Cells.Item($srcRow, $col).Copy($ws.Cells.Item($newRow, $col))
}

$ws.Cells.Item($newRow, 1).Value = "2025-08-26 06:50:16 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-26 12:20:16 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
